{"js": "// Replace each math-problem cell's text with its updated value, in\n// document order (row-major), preserving all existing run formatting\n// (font, size, etc.) by writing through Table.values.\nconst newValues = [\n  [\"96-11=85\", \"96-68=28\", \"39-12=27\", \"27+41=68\", \"69-61=8\"],\n  [\"85-50=35\", \"53-7=46\", \"76-32=44\", \"38+20=58\", \"34+52=86\"],\n  [\"83-45=38\", \"61+37=98\", \"54-46=8\", \"85-19=66\", \"48-32=16\"],\n  [\"38+1=39\", \"89-87=2\", \"11+17=28\", \"79-78=1\", \"28+3=31\"],\n  [\"9-7=2\", \"81-34=47\", \"61-12=49\", \"24+62=86\", \"83+3=86\"],\n  [\"22+14=36\", \"29+18=47\", \"86-53=33\", \"25+26=51\", \"92-76=16\"],\n  [\"77+5=82\", \"67-43=24\", \"65+26=91\", \"14+37=51\", \"6+11=17\"],\n  [\"54+26=80\", \"68-43=25\", \"34+24=58\", \"47+49=96\", \"18+29=47\"],\n  [\"62-49=13\", \"96-74=22\", \"23+23=46\", \"50-27=23\", \"59+8=67\"],\n  [\"50+35=85\", \"47+17=64\", \"3+55=58\", \"84-83=1\", \"9+31=40\"],\n  [\"92-63=29\", \"93-82=11\", \"83+5=88\", \"12-8=4\", \"37-30=7\"],\n  [\"43+31=74\", \"22+26=48\", \"71-24=47\", \"19+34=53\", \"63-20=43\"],\n  [\"36+52=88\", \"48+38=86\", \"58+27=85\", \"56+38=94\", \"95-32=63\"],\n  [\"27-14=13\", \"59+21=80\", \"67-13=54\", \"62-23=39\", \"6+19=25\"],\n  [\"4+39=43\", \"52+10=62\", \"70+25=95\", \"52+35=87\", \"30+8=38\"],\n  [\"75-36=39\", \"73-12=61\", \"5+44=49\", \"30+17=47\", \"12+60=72\"],\n  [\"50-27=23\", \"30+1=31\", \"15-11=4\", \"59-12=47\", \"43+33=76\"],\n  [\"3+90=93\", \"50-11=39\", \"81+16=97\", \"88-50=38\", \"60-14=46\"],\n  [\"0+40=40\", \"81+5=86\", \"76-27=49\", \"71-13=58\", \"2+85=87\"],\n  [\"75-17=58\", \"55-3=52\", \"28-15=13\", \"16+30=46\", \"75-34=41\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace each math-problem cell's text with its updated value, in\n# document order (row-major), preserving run formatting by writing\n# through Cell.Range.Text (in-place text replacement).\n$newValues = @(\n    @('96-11=85', '96-68=28', '39-12=27', '27+41=68', '69-61=8'),\n    @('85-50=35', '53-7=46', '76-32=44', '38+20=58', '34+52=86'),\n    @('83-45=38', '61+37=98', '54-46=8', '85-19=66', '48-32=16'),\n    @('38+1=39', '89-87=2', '11+17=28', '79-78=1', '28+3=31'),\n    @('9-7=2', '81-34=47', '61-12=49', '24+62=86', '83+3=86'),\n    @('22+14=36', '29+18=47', '86-53=33', '25+26=51', '92-76=16'),\n    @('77+5=82', '67-43=24', '65+26=91', '14+37=51', '6+11=17'),\n    @('54+26=80', '68-43=25', '34+24=58', '47+49=96', '18+29=47'),\n    @('62-49=13', '96-74=22', '23+23=46', '50-27=23', '59+8=67'),\n    @('50+35=85', '47+17=64', '3+55=58', '84-83=1', '9+31=40'),\n    @('92-63=29', '93-82=11', '83+5=88', '12-8=4', '37-30=7'),\n    @('43+31=74', '22+26=48', '71-24=47', '19+34=53', '63-20=43'),\n    @('36+52=88', '48+38=86', '58+27=85', '56+38=94', '95-32=63'),\n    @('27-14=13', '59+21=80', '67-13=54', '62-23=39', '6+19=25'),\n    @('4+39=43', '52+10=62', '70+25=95', '52+35=87', '30+8=38'),\n    @('75-36=39', '73-12=61', '5+44=49', '30+17=47', '12+60=72'),\n    @('50-27=23', '30+1=31', '15-11=4', '59-12=47', '43+33=76'),\n    @('3+90=93', '50-11=39', '81+16=97', '88-50=38', '60-14=46'),\n    @('0+40=40', '81+5=86', '76-27=49', '71-13=58', '2+85=87'),\n    @('75-17=58', '55-3=52', '28-15=13', '16+30=46', '75-34=41')\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n\n"}
